# Applies the "below_ground.xlsx" param_dist update:
#  - adds new rows (corrosivity, install_year, PG&E liquefaction model params,
#    repair-rate model params, etc.) to the "fixed" sheet
#  - moves the active tab/selection from "level1" to "fixed"

$wb = $excel.ActiveWorkbook
$wsFixed = $wb.Worksheets.Item("fixed")
$wsLevel1 = $wb.Worksheets.Item("level1")

# ---------------------------------------------------------------------------
# 1. Populate the new rows (7-16) on the "fixed" sheet.
#    Cell values are written in the same order the strings were first typed
#    in the original edit session so that new shared-string entries land in
#    the expected sequence.
# ---------------------------------------------------------------------------

# Row 7 - corrosivity
$wsFixed.Range("A7").Value = "corrosivity"
$wsFixed.Range("D7").Value = "soil corrosivity: <=750 = high impact, >1500 = low impact"
$wsFixed.Range("E7").Value = "ohm-cm"
$wsFixed.Range("B7").Value = $true
$wsFixed.Range("C7").Value = $false
$wsFixed.Range("F7").Value = 1600

# Row 8 - install_year
$wsFixed.Range("A8").Value = "install_year"
$wsFixed.Range("D8").Value = "year of installation"
$wsFixed.Range("E8").Value = "unitless"
$wsFixed.Range("B8").Value = $true
$wsFixed.Range("C8").Value = $false
$wsFixed.Range("F8").Value = 1980

# Row 9 - pge_a
$wsFixed.Range("A9").Value = "pge_a"
$wsFixed.Range("D9").Value = "PG&E model coefficent a"
$wsFixed.Range("E9").Value = "probability (fraction)"
$wsFixed.Range("B9").Value = $true
$wsFixed.Range("C9").Value = $false
$wsFixed.Range("F9").Value = "internal gis dataset"

# Column A for rows 10 and 11 (pge_b, pge_c) were typed before their
# descriptions in column D.
$wsFixed.Range("A10").Value = "pge_b"
$wsFixed.Range("A11").Value = "pge_c"

$wsFixed.Range("D10").Value = "PG&E model coefficent b"
$wsFixed.Range("D11").Value = "PG&E model coefficent c"

# Descriptions for rows 12, 13, 15, 14, 16 (note the 15/14 order) were typed
# next, followed by the rv_label column A values for rows 12-16.
$wsFixed.Range("D12").Value = "PG&E model maximum magnitude scaling factor"
$wsFixed.Range("D13").Value = "PG&E model lateral-spreading coefficient"
$wsFixed.Range("D15").Value = "PG&E model settlement coefficent"
$wsFixed.Range("D14").Value = "PG&E model lateral-spreading sigmaln"
$wsFixed.Range("D16").Value = "PG&E model settlement sigmaln"

$wsFixed.Range("A12").Value = "msf_max"
$wsFixed.Range("A13").Value = "dl"
$wsFixed.Range("A14").Value = "sigl"
$wsFixed.Range("A15").Value = "ds"
$wsFixed.Range("A16").Value = "sigs"

# Last new shared string: the "g" unit for row 10.
$wsFixed.Range("E10").Value = "g"

# Remaining (reused) values for rows 10-16.
$wsFixed.Range("B10").Value = $true
$wsFixed.Range("C10").Value = $false
$wsFixed.Range("F10").Value = "internal gis dataset"

$wsFixed.Range("B11").Value = $true
$wsFixed.Range("C11").Value = $false
$wsFixed.Range("E11").Value = "unitless"
$wsFixed.Range("F11").Value = "internal gis dataset"

$wsFixed.Range("B12").Value = $true
$wsFixed.Range("C12").Value = $false
$wsFixed.Range("E12").Value = "unitless"
$wsFixed.Range("F12").Value = "internal gis dataset"

$wsFixed.Range("B13").Value = $true
$wsFixed.Range("C13").Value = $false
$wsFixed.Range("E13").Value = "m"
$wsFixed.Range("F13").Value = "internal gis dataset"

$wsFixed.Range("B14").Value = $true
$wsFixed.Range("C14").Value = $false
$wsFixed.Range("E14").Value = "unitless"
$wsFixed.Range("F14").Value = "internal gis dataset"

$wsFixed.Range("B15").Value = $true
$wsFixed.Range("C15").Value = $false
$wsFixed.Range("E15").Value = "m"
$wsFixed.Range("F15").Value = "internal gis dataset"

$wsFixed.Range("B16").Value = $true
$wsFixed.Range("C16").Value = $false
$wsFixed.Range("E16").Value = "unitless"
$wsFixed.Range("F16").Value = "internal gis dataset"

# ---------------------------------------------------------------------------
# 2. Update the active sheet / selection so that "fixed" becomes the
#    selected tab (previously "level1" was selected).
# ---------------------------------------------------------------------------

$wsLevel1.Activate() | Out-Null
$wsLevel1.Range("F20").Select() | Out-Null

$wsFixed.Activate() | Out-Null
$wsFixed.Range("F15").Select() | Out-Null
